$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Obrigatorio" column (E) from "N" to "S" for rows 2 through 9
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 5).Value = "S"
}
